$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 155, pushing old rows 155-168 down to 156-169.
$ws.Rows("155:155").Insert() | Out-Null

# The new row 155 gets the data that used to be in row 154 (Lauren Moran).
$ws.Range("A155").Value = "Lauren"
$ws.Range("B155").Value = "Moran"
$ws.Range("C155").Value = "Maywood"
$ws.Range("D155").Value = "Maywood, New Jersey"
$ws.Range("E155").Value = "USA"
$ws.Range("H155").Value = "Lauren Moran.jpg"

# Row 154 becomes the new roster entry: Kaitlin Moran.
$ws.Range("A154").Value = "Kaitlin"
$ws.Range("H154").Value = "Kaitlin Moran.jpg"

# Restore the view state to match the edited workbook.
$ws.Application.ActiveWindow.ScrollRow = 150
$ws.Range("K156").Select() | Out-Null
